# Add team record (Wins/Losses/Ties) columns to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in AD1:AF1, matching the style used by the existing header row (A1 etc.)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the win/loss/tie record for each data row (2 through 47)
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 69   # AD
    $ws.Cells.Item($row, 31).Value = 93   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
